$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.245.71'
$ws.Range('E2').Value = '  -5.68%  '
$ws.Range('D3').Value = '2.458.04'
$ws.Range('E3').Value = '  -8.22%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '540.10'
$ws.Range('E5').Value = '  -2.48%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '146.23'
$ws.Range('E6').Value = '  -6.95%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('E8').Value = '  -2.18%  '
$ws.Range('D9').Value = '2.474.76'
$ws.Range('E9').Value = '  -7.83%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0999'
$ws.Range('E10').Value = '  -4.97%  '
$ws.Range('E11').Value = '  -1.26%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.51'
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.353'
$ws.Range('E13').Value = '  -3.55%  '
$ws.Range('D14').Value = '2.895.57'
$ws.Range('E14').Value = '  -8.24%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '24.27'
$ws.Range('E15').Value = '  -7.15%  '
$ws.Range('D16').Value = '59.224.99'
$ws.Range('E16').Value = '  -5.60%  '
$ws.Range('E17').Value = '  -4.87%  '
$ws.Range('D18').Value = '2.469.77'
$ws.Range('E18').Value = '  -7.91%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.22'
$ws.Range('E19').Value = '  -4.96%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.38'
$ws.Range('E20').Value = '  -4.40%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '323.95'
$ws.Range('E21').Value = '  -5.78%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.982'
$ws.Range('E22').Value = '  -1.83%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.76'
$ws.Range('E23').Value = '  -7.04%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '60.90'
$ws.Range('E24').Value = '  -3.68%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.453'
$ws.Range('E25').Value = '  -10.82%  '
$ws.Range('E26').Value = '  -5.01%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.989'
$ws.Range('E27').Value = '  -0.97%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.76'
$ws.Range('E28').Value = '  -4.63%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.80'
$ws.Range('E29').Value = '  -6.22%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.27'
$ws.Range('E30').Value = '  -9.12%  '
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').Value = '0.0₃0778'
$ws.Range('E31').Value = '  -8.36%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.81'
$ws.Range('E32').Value = '  -5.78%  '
$ws.Range('E33').Value = '  -0.18%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '157.94'
$ws.Range('E34').Value = '  -2.98%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.41'
$ws.Range('E35').Value = '  -2.44%  '
$ws.Range('E36').Value = '  -3.73%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.47'
$ws.Range('E37').Value = '  -7.87%  '
$ws.Range('E38').Value = '  -3.28%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.88'
$ws.Range('E39').Value = '  -4.10%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '313.44'
$ws.Range('E40').Value = '  -7.33%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '36.46'
$ws.Range('E41').Value = '  -4.85%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.74'
$ws.Range('E42').Value = '  -6.16%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.831'
$ws.Range('E43').Value = '  -9.83%  '
$ws.Range('E44').Value = '  -0.46%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.598'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '10.74'
$ws.Range('E46').Value = '  -2.31%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '125.92'
$ws.Range('E47').Value = '  -3.22%  '
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0530'
$ws.Range('E48').Value = '  -4.46%  '
$ws.Range('E49').Value = '  -3.33%  '
$ws.Range('E50').Value = '  -3.64%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '18.53'
$ws.Range('E51').Value = '  -7.98%  '
